$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5263440
$ws.Range("I6").Value = 6579100
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 19737300
$ws.Range("L6").Value = 2400
$ws.Range("M6").Value = -19737188
$ws.Range("N6").Value = -2624

$ws.Range("H9").Value = 80.666664
$ws.Range("I9").Value = 45.333332
$ws.Range("J9").Value = 98.333336
$ws.Range("K9").Value = 45.333332
$ws.Range("L9").Value = 98.333336
$ws.Range("M9").Value = 123.666668
$ws.Range("N9").Value = -436.333336

$ws.Range("H69").Value = 20002388
$ws.Range("I69").Value = 2466.6667
$ws.Range("J69").Value = 21278978
$ws.Range("K69").Value = 7400.000100000001
$ws.Range("L69").Value = 63836934
$ws.Range("M69").Value = -6526.000100000001
$ws.Range("N69").Value = -63838682

$ws.Range("H72").Value = 20002388
$ws.Range("I72").Value = 2466.6667
$ws.Range("J72").Value = 21278978
$ws.Range("K72").Value = 22200.0003
$ws.Range("L72").Value = 191510802
$ws.Range("M72").Value = -17832.0003
$ws.Range("N72").Value = -191519538

$ws.Range("H141").Value = 1092.7273
$ws.Range("I141").Value = 1103
$ws.Range("J141").Value = 990
$ws.Range("K141").Value = 3309
$ws.Range("L141").Value = 2970
$ws.Range("M141").Value = 1871
$ws.Range("N141").Value = -13330

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 15000
$ws.Range("J76").Value = 15000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15676

$ws.Range("H79").Value = 15000
$ws.Range("J79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17340

$ws.Range("H102").Value = 2576.375
$ws.Range("I102").Value = 3052.75
$ws.Range("J102").Value = 2100
$ws.Range("K102").Value = 3052.75
$ws.Range("L102").Value = 2100
$ws.Range("M102").Value = -1430.75
$ws.Range("N102").Value = -5344

$ws.Range("H132").Value = 2830.6191
$ws.Range("I132").Value = 1322.5555
$ws.Range("J132").Value = 3961.6667
$ws.Range("K132").Value = 3967.6665
$ws.Range("L132").Value = 11885.0001
$ws.Range("M132").Value = -1437.6665
$ws.Range("N132").Value = -16945.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1869.9286
$ws.Range("I86").Value = 1953.3334
$ws.Range("J86").Value = 1719.8
$ws.Range("K86").Value = 1953.3334
$ws.Range("L86").Value = 1719.8
$ws.Range("M86").Value = -830.3334
$ws.Range("N86").Value = -3965.8

$ws.Range("H89").Value = 1869.9286
$ws.Range("I89").Value = 1953.3334
$ws.Range("J89").Value = 1719.8
$ws.Range("K89").Value = 9766.666999999999
$ws.Range("L89").Value = 8599
$ws.Range("M89").Value = -4150.666999999999
$ws.Range("N89").Value = -19831

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2732.8276
$ws.Range("I31").Value = 1550.3715
$ws.Range("J31").Value = 4532.2173
$ws.Range("K31").Value = 1550.3715
$ws.Range("L31").Value = 4532.2173
$ws.Range("M31").Value = -1255.3715
$ws.Range("N31").Value = -5122.2173

$ws.Range("H34").Value = 2732.8276
$ws.Range("I34").Value = 1550.3715
$ws.Range("J34").Value = 4532.2173
$ws.Range("K34").Value = 1550.3715
$ws.Range("L34").Value = 4532.2173
$ws.Range("M34").Value = -1348.3715
$ws.Range("N34").Value = -4936.2173

$ws.Range("H62").Value = 2592.8572
$ws.Range("I62").Value = 2412.5
$ws.Range("J62").Value = 2833.3333
$ws.Range("K62").Value = 2412.5
$ws.Range("L62").Value = 2833.3333
$ws.Range("M62").Value = -1788.5
$ws.Range("N62").Value = -4081.3333

$ws.Range("H65").Value = 2592.8572
$ws.Range("I65").Value = 2412.5
$ws.Range("J65").Value = 2833.3333
$ws.Range("K65").Value = 12062.5
$ws.Range("L65").Value = 14166.6665
$ws.Range("M65").Value = -8942.5
$ws.Range("N65").Value = -20406.6665

$ws.Range("H134").Value = 3172.6
$ws.Range("I134").Value = 1620.1538
$ws.Range("J134").Value = 4854.4165
$ws.Range("K134").Value = 4860.4614
$ws.Range("L134").Value = 14563.2495
$ws.Range("M134").Value = -2325.4614
$ws.Range("N134").Value = -19633.2495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 706.9722
$ws.Range("J5").Value = 1316.6666
$ws.Range("L5").Value = 3949.9998
$ws.Range("N5").Value = -4173.9998

$ws.Range("H59").Value = 1223
$ws.Range("I59").Value = 371.66666
$ws.Range("K59").Value = 1114.99998
$ws.Range("M59").Value = -574.9999800000001

$ws.Range("H80").Value = 2042.8572
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2042.8572
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 6128.571599999999
$ws.Range("N80").Value = -8000.571599999999

$ws.Range("H83").Value = 2042.8572
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2042.8572
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 18385.7148
$ws.Range("N83").Value = -27745.7148

$ws.Range("H92").Value = 657.5714
$ws.Range("J92").Value = 657.5714
$ws.Range("L92").Value = 1972.7142
$ws.Range("N92").Value = -4468.7142

$ws.Range("H110").Value = 4050.25

$ws.Range("H112").Value = 2556.9092
$ws.Range("J112").Value = 2900
$ws.Range("L112").Value = 8700
$ws.Range("N112").Value = -10916

$ws.Range("H113").Value = 3831973.2
$ws.Range("I113").Value = 8621065
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 25863195
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = -25861025
$ws.Range("N113").Value = -6440

$ws.Range("H115").Value = 632.1429000000001
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H117").Value = 2666.2856
$ws.Range("I117").Value = 2527.3333
$ws.Range("J117").Value = 3500
$ws.Range("K117").Value = 7581.999899999999
$ws.Range("L117").Value = 10500
$ws.Range("M117").Value = -4139.999899999999
$ws.Range("N117").Value = -17384

$ws.Range("H118").Value = 1985.8
$ws.Range("I118").Value = 876.3333
$ws.Range("J118").Value = 3650
$ws.Range("K118").Value = 2628.9999
$ws.Range("L118").Value = 10950
$ws.Range("M118").Value = -1385.9999
$ws.Range("N118").Value = -13436

$ws.Range("H121").Value = 1033.1
$ws.Range("I121").Value = 647.5
$ws.Range("J121").Value = 1173.3182
$ws.Range("K121").Value = 1942.5
$ws.Range("L121").Value = 3519.9546
$ws.Range("M121").Value = -632.5
$ws.Range("N121").Value = -6139.9546

$ws.Range("H125").Value = 3462
$ws.Range("I125").Value = 1924
$ws.Range("J125").Value = 5000
$ws.Range("K125").Value = 5772
$ws.Range("L125").Value = 15000
$ws.Range("M125").Value = -852
$ws.Range("N125").Value = -24840

$ws.Range("H135").Value = 706.9722
$ws.Range("J135").Value = 1316.6666
$ws.Range("L135").Value = 11849.9994
$ws.Range("N135").Value = -16919.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2915.2856
$ws.Range("I80").Value = 2723.75
$ws.Range("J80").Value = 3033.1538
$ws.Range("K80").Value = 2723.75
$ws.Range("L80").Value = 3033.1538
$ws.Range("M80").Value = -1725.75
$ws.Range("N80").Value = -5029.1538

$ws.Range("H83").Value = 2915.2856
$ws.Range("I83").Value = 2723.75
$ws.Range("J83").Value = 3033.1538
$ws.Range("K83").Value = 13618.75
$ws.Range("L83").Value = 15165.769
$ws.Range("M83").Value = -8626.75
$ws.Range("N83").Value = -25149.769

$ws.Range("H113").Value = 931791.9
$ws.Range("I113").Value = 1558.7142
$ws.Range("J113").Value = 1582955.1
$ws.Range("K113").Value = 1558.7142
$ws.Range("L113").Value = 1582955.1
$ws.Range("M113").Value = 611.2858000000001
$ws.Range("N113").Value = -1587295.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2022.4814
$ws.Range("I68").Value = 1965.35
$ws.Range("J68").Value = 2185.7144
$ws.Range("K68").Value = 1965.35
$ws.Range("L68").Value = 2185.7144
$ws.Range("M68").Value = -1216.35
$ws.Range("N68").Value = -3683.7144

$ws.Range("H71").Value = 2022.4814
$ws.Range("I71").Value = 1965.35
$ws.Range("J71").Value = 2185.7144
$ws.Range("K71").Value = 9826.75
$ws.Range("L71").Value = 10928.572
$ws.Range("M71").Value = -6082.75
$ws.Range("N71").Value = -18416.572

$ws.Range("H82").Value = 1599.4445
$ws.Range("I82").Value = 1503.88
$ws.Range("J82").Value = 1816.6364
$ws.Range("K82").Value = 1503.88
$ws.Range("L82").Value = 1816.6364
$ws.Range("M82").Value = -1142.88
$ws.Range("N82").Value = -2538.6364

$ws.Range("H85").Value = 1599.4445
$ws.Range("I85").Value = 1503.88
$ws.Range("J85").Value = 1816.6364
$ws.Range("K85").Value = 1503.88
$ws.Range("L85").Value = 1816.6364
$ws.Range("M85").Value = -255.8800000000001
$ws.Range("N85").Value = -4312.6364

$ws.Range("H122").Value = 2011.9131
$ws.Range("I122").Value = 1859.625
$ws.Range("J122").Value = 2360
$ws.Range("K122").Value = 5578.875
$ws.Range("L122").Value = 7080
$ws.Range("M122").Value = -3128.875
$ws.Range("N122").Value = -11980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2417.4443
$ws.Range("I132").Value = 1210.2727
$ws.Range("J132").Value = 4314.4287
$ws.Range("K132").Value = 3630.8181
$ws.Range("L132").Value = 12943.2861
$ws.Range("M132").Value = -1100.8181
$ws.Range("N132").Value = -18003.2861
